$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled/recalculated data
$ws.Range("F2").Value = 7
$ws.Range("F5").Value = -8
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = -5
$ws.Range("F8").Value = -8
$ws.Range("F9").Value = -5
$ws.Range("F10").Value = 2
$ws.Range("F14").Value = 8
$ws.Range("F16").Value = 2
